# varInfo and createStructure interplay
#
# Inserts a new "Unterteilung.im.Skalenhandbuch" column after
# "in.DS.und.SH" (shifting Layout..Seitenumbruch.im.Inhaltsverzeichnis one
# column to the right), fills it in for the 3 data rows, tweaks the
# Anmerkung.Var value for VAR1, and turns the Gliederung column from a
# plain row-order number into a text "x.xx" code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("Unterteilung.im.Skalenhandbuch"). Everything
# that used to live in C..P (Layout .. Seitenumbruch.im.Inhaltsverzeichnis)
# shifts right into D..Q.
$ws.Columns("C").Insert()

# Header for the newly inserted column
$ws.Range("C1").Value = "Unterteilung.im.Skalenhandbuch"

# Values for the new column
$ws.Range("C2").Value = "Teil 1a"
$ws.Range("C3").Value = "Teil 1b"
$ws.Range("C4").Value = "Teil 2"

# Anmerkung.Var (now column F) - VAR1's row changes from "NA" to "-"
$ws.Range("F2").Value = "-"

# Gliederung (now column G) switches from a numeric layout order (1/2/3)
# to a textual "x.xx" code. Force text storage (instead of Excel's default
# numeric auto-detection) by toggling the number format around the write,
# then clearing the cell format again so no formatting change is left
# behind.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1.01"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "1.02"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "2.01"

$ws.Range("G2:G4").ClearFormats()
